$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell B11 used to hold the shared string "R40". It now holds the literal
# text "1" (a new shared-string entry). Prefixing with an apostrophe tells
# Excel to store the numeric-looking value as text (t="s") rather than
# coercing it to a number, exactly like typing '1 into the cell.
$ws.Range("B11").Value = "'1"
